$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date values in column A (keep the existing date number format / style)
$ws.Range("A2").Value = 46062
$ws.Range("A5").Value = 46066
$ws.Range("A9").Value = 46069
$ws.Range("A13").Value = 46076
$ws.Range("A17").Value = 46101
$ws.Range("A20").Value = 46121

# Update the corresponding weekday labels in column B to match the new dates
$ws.Range("B2").Value = "lundi"
$ws.Range("B5").Value = "vendredi"
$ws.Range("B9").Value = "lundi"
$ws.Range("B13").Value = "lundi"
$ws.Range("B17").Value = "vendredi"
$ws.Range("B20").Value = "jeudi"
